$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Strip the old uniform cell style ("style 2") from content + from the
# column definitions (A:D) so the rebuilt grid uses plain default formatting,
# matching the target file (only A1 keeps an explicit style, and it becomes
# a quotePrefix style once we enter literal text starting with "=").
$ws.Cells.Clear()
$ws.Columns("A:D").ClearFormats()

# --- Structural edits: a new blank row on top, and a new "Flow" column
# inserted before the old last (Test/Test case) column.
$ws.Rows("1:1").Insert()
$ws.Columns("D:D").Insert()

# New column D ("Flow") should carry the same width as column C (18.77734375
# raw units); the COM width setter only has ~1/6-character resolution, so we
# use the closest reproducible value.
$ws.Columns("D:D").ColumnWidth = 18

# --- Row 1: a literal note that happens to start with "=", entered as text
# (leading apostrophe forces literal/quoted text instead of a formula).
$ws.Range("A1").Value = "'=> Ausführung eines Testplanes noch nicht implementiert - wird via listener die Testcases on the fly erstellen und dann ausführen"

# --- Row 2: headers
$ws.Range("A2").Value = "Run Y/N?"
$ws.Range("B2").Value = "Test area"
$ws.Range("C2").Value = "Test object"
$ws.Range("D2").Value = "Flow"
$ws.Range("E2").Value = "Test case"

# --- Row 3
$ws.Range("A3").Value = "Y"
$ws.Range("B3").Value = "001_Login"
$ws.Range("C3").Value = "001_Login"
$ws.Range("D3").Value = "floLogin"
$ws.Range("E3").Value = "001_Login_001_Successful"

# --- Row 4
$ws.Range("A4").Value = "N"
$ws.Range("D4").Value = "floLogin"
$ws.Range("E4").Value = "001_Login_002_LoginWithoutEmailAndPassword"

# --- Row 5 (B5 keeps the distinct placeholder look previously on B4)
$ws.Range("A5").Value = "N"
$ws.Range("B5").Font.Name = "Inter"
$ws.Range("B5").Font.Size = 10
$ws.Range("B5").Font.Color = 12893884
$ws.Range("B5").VerticalAlignment = -4108
$ws.Range("C5").Value = "002_PasswordReset"
$ws.Range("D5").Value = "floLogin"
$ws.Range("E5").Value = "002_PasswordReset_001_Successful"

# --- Row 6
$ws.Range("A6").Value = "N"
$ws.Range("D6").Value = "floLogin"
$ws.Range("E6").Value = "002_PasswordReset_002_BackToLogin"

# --- Row 7
$ws.Range("A7").Value = "N"
$ws.Range("B7").Value = "003_MeinProfil"
$ws.Range("C7").Value = "001_Allgemein"
$ws.Range("D7").Value = "floProfil"
$ws.Range("E7").Value = "003_Profil_001_Allgemein_CheckingDefaults"

# --- Row 8 (new)
$ws.Range("A8").Value = "N"
$ws.Range("C8").Value = "002_Profil"
$ws.Range("D8").Value = "floProfil"
$ws.Range("E8").Value = "003_Profil_002_Profil_Datenaenderung"

# --- Row 9 (new)
$ws.Range("A9").Value = "N"
$ws.Range("C9").Value = "003_Abwesenheiten"
$ws.Range("D9").Value = "floProfil"
$ws.Range("E9").Value = "003_Profil_003_Abwesenheiten_Normalfall_Anlage"

# --- Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection
$ws.Range("E19").Select() | Out-Null
